# Append a newly-scraped Lancers listing (run timestamp 2026-02-09 02:00:33).
#
# The sheet is sorted by "優先度スコア" (priority score), and the new
# listing (score 38) ranks 3rd, so it is inserted as row 5 - pushing the
# previous rows 5 ("bubble...") and 6 ("BigQuery+Looker Studio...") down to
# rows 6 and 7. Every data row's "取得日時" (fetched-at) timestamp is
# refreshed to the new run time as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-02-09 02:00:33"

# 1. Insert a new row at row 5; this shifts the old rows 5-6 down to 6-7,
#    carrying their values/styles/hyperlink anchors along with them.
$ws.Rows("5:5").Insert()

# 2. Populate the newly inserted row 5 with the new listing's data.
$ws.Range("A5").Value = $newTimestamp
$ws.Range("B5").Value = "養鰻管理Excelの判断ロジック(給餌)を理解し、継続的に伴走できる方を募集"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5488109"
$ws.Range("G5").Value = 38
$ws.Range("H5").Value = "◇管理"

# 3. Refresh "取得日時" on every other existing data row to the new run time.
$ws.Range("A2").Value = $newTimestamp
$ws.Range("A3").Value = $newTimestamp
$ws.Range("A4").Value = $newTimestamp
$ws.Range("A6").Value = $newTimestamp
$ws.Range("A7").Value = $newTimestamp

# 4. Rebuild the hyperlinks collection so every F-cell's relationship
#    target matches the URL now displayed in that cell (Insert() only
#    shifts a hyperlink's cell anchor, it does not retarget it, and the
#    row shift means rows 5 and 6 now show different URLs than before).
$hlStyle = $ws.Range("F2").Style
$ws.Hyperlinks.Delete()

$urls = @(
    "https://www.lancers.jp/work/detail/5487791",
    "https://www.lancers.jp/work/detail/5487945",
    "https://www.lancers.jp/work/detail/5487838",
    "https://www.lancers.jp/work/detail/5488109",
    "https://www.lancers.jp/work/detail/5487908",
    "https://www.lancers.jp/work/detail/5487828"
)
for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = $i + 2
    $ws.Hyperlinks.Add($ws.Range("F$row"), $urls[$i]) | Out-Null
}
$ws.Range("F2:F7").Style = $hlStyle

$wb.Save()
